$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.680.74'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '3.786.81'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("D7").Value = '3.773.21'
$ws.Range("E7").Value = '  +0.81%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.92%  '
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000252'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").Value = '4.423.49'
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("D16").Value = '3.793.04'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.89%  '
$ws.Range("D18").Value = '67.650.04'
$ws.Range("E18").Value = '  -0.95%  '
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.696'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("E24").Value = '  +4.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.68%  '
$ws.Range("E26").Value = '  +0.95%  '
$ws.Range("E27").Value = '  -2.71%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D30").Value = '3.932.93'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.77'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.65%  '
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("E34").Value = '  -0.69%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.09'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.995'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.95%  '
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '393.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.73'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.34%  '
